$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.971.51"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.912.68"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'355.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("D7").Value = "'0.569"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.45%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.629"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "'38.83"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("D11").Value = "'0.138"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "'0.0869"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "'19.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "'7.77"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "3.366.91"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "2.900.95"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "'0.983"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "51.925.73"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'3.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.28%  "
$ws.Range("D20").Value = "'7.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "'13.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "'70.46"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "'268.37"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "'2.81"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("E26").Value = "  +9.24%  "
$ws.Range("D27").Value = "'26.86"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("D28").Value = "'7.64"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +16.49%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'0.106"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +8.53%  "
$ws.Range("D31").Value = "'10.48"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'37.52"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "'6.16"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").Value = "'52.17"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").Value = "'18.18"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").Value = "'1.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("D41").Value = "'2.71"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.56%  "
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("D43").Value = "'22.96"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").Value = "'119.50"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -6.55%  "
$ws.Range("D47").Value = "'3.45"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "2.126.94"
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("D49").Value = "'0.253"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.84%  "
$ws.Range("D50").Value = "'0.0337"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "'0.925"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.57%  "
